$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Refresh STIX IDs (column B, rows 2-83) to the new attack-pattern UUIDs ---
#     (ATLAS Data 4.5.1 / STIX regenerated against ATT&CK v14.1)
$newStixIds = @(
    "attack-pattern--4fb3f968-42a6-43da-ae93-5f165afa82e3",
    "attack-pattern--99ee109d-c2d0-48cc-9960-d1d2a9ecd0ad",
    "attack-pattern--68b13c6b-c45d-4864-8186-189c77eaeb4d",
    "attack-pattern--4d881f79-6a16-4e84-a795-4f3914045e67",
    "attack-pattern--b6a3ddcc-60d7-4ae8-b702-eca52e3ace5d",
    "attack-pattern--f4242e1a-203a-4b26-8a7f-967458b3751a",
    "attack-pattern--1a328de1-1327-4ed3-98c0-0a2e9da558c0",
    "attack-pattern--0dc3c131-1524-486e-a653-52240e95a22a",
    "attack-pattern--e641274e-cac2-49c4-800f-abd0c0d825cd",
    "attack-pattern--047c1235-9c18-4425-8494-6f8a75feabce",
    "attack-pattern--35bd4e6e-9056-4115-9395-81735cfd1761",
    "attack-pattern--d93718ee-12c5-4f76-9845-1bb14f948513",
    "attack-pattern--36a77758-5a44-48e9-9545-14bd0a54756d",
    "attack-pattern--acebf60b-3f22-4c75-a74f-4da8788fae17",
    "attack-pattern--ff90c91c-2f51-48db-bc3c-e1d1ebf093f3",
    "attack-pattern--2e3499bf-169b-4b77-ae04-c689e53a9895",
    "attack-pattern--0ba8d50c-bbcc-4b6b-ac2d-0d78fdae0422",
    "attack-pattern--50d26291-d3ac-451f-b051-3159ba12f39e",
    "attack-pattern--05c6f7a3-1f8a-443d-a3d7-b4d54ab28e8c",
    "attack-pattern--87f481a5-f4f5-4aa8-b19a-854634bdf08f",
    "attack-pattern--61527318-c59b-4097-840c-adf9961bad42",
    "attack-pattern--c3e09f44-3ff9-4db9-bf26-6a42bbcfbbb1",
    "attack-pattern--a5f36bcf-f003-4bde-bd38-873ecb1160c6",
    "attack-pattern--501caefa-6ffd-496b-be90-c644e9dc8037",
    "attack-pattern--c3c55a49-99c1-4613-994a-9948228ef0c8",
    "attack-pattern--dbd25a74-6024-4e30-9ba2-20428a447b70",
    "attack-pattern--a0205c5e-9527-4020-8927-26c855a45008",
    "attack-pattern--37d808a2-ca53-4b8f-a5c9-1dfc4bc59d6f",
    "attack-pattern--5c4401a9-b406-4cbd-9c2d-af905a927c23",
    "attack-pattern--22fd9008-0257-448b-a013-e253ece129e7",
    "attack-pattern--26e0f122-c1ef-4d64-b20a-5cd81b1a1979",
    "attack-pattern--ef5867e1-c9cc-4e52-8044-80053e308042",
    "attack-pattern--8e2247ab-ddb9-4105-9bca-7dc756987394",
    "attack-pattern--4e33fd0b-b624-4aa9-9c87-69acaa29437d",
    "attack-pattern--7adc32b4-93f3-441f-9ad4-08f8643d450e",
    "attack-pattern--1b7c868f-f99c-4b9c-98c2-fe8608cc39c4",
    "attack-pattern--0c25c07a-0800-4635-ac9c-fa3fbd0b571d",
    "attack-pattern--36d540cd-e31b-457e-8ac0-f34eec004a7a",
    "attack-pattern--4af77ce4-5699-4bc0-85e5-9a7b115a9fda",
    "attack-pattern--8eb6f82a-6c90-44a8-bd8d-a0ba632858a5",
    "attack-pattern--d59358df-a685-4a5d-b65b-2f81bd91bc02",
    "attack-pattern--80a3250a-c57b-445e-a668-c9bb8bf1b136",
    "attack-pattern--657cba2d-a94b-46ed-96c6-b0ff825e79d6",
    "attack-pattern--e04885b8-82ed-464c-9789-fafb7fdbbe30",
    "attack-pattern--db789568-054a-498b-b585-56e8d5ce88f9",
    "attack-pattern--a0da3f30-9270-418f-834e-7c79bed1065e",
    "attack-pattern--ccd395a9-69f3-4f5e-a819-99f98bcc1c4b",
    "attack-pattern--69268719-393f-41d3-801e-18d607474b01",
    "attack-pattern--32287d8a-ac22-4d9f-a4ba-a55566fa8b36",
    "attack-pattern--2308375d-743a-49ea-9fba-509a14e4a5e2",
    "attack-pattern--21a73c23-60b5-46bc-b480-afa101fbac5f",
    "attack-pattern--7344ac61-c0ef-484b-acb3-8ae385fcddb8",
    "attack-pattern--4e479550-b51e-4399-bdda-4895d98f99fd",
    "attack-pattern--7f8b61c0-e4a2-4a21-851a-7f8041a929c6",
    "attack-pattern--2fb39d9e-2363-4672-a3c8-4c6c63c6b94c",
    "attack-pattern--451cf23a-2346-4c3c-aa15-a629af8bb982",
    "attack-pattern--aab01876-c6ba-4fac-a23a-25bc508e6a5c",
    "attack-pattern--12b09c3e-46a9-4647-9df3-68cfde54ba21",
    "attack-pattern--9d8044ed-6082-44d8-bc62-559f9296b580",
    "attack-pattern--2f154be7-a9fa-4c5c-b446-5c4ca11049d7",
    "attack-pattern--4719e4eb-78c9-4292-8059-a88886b65dfe",
    "attack-pattern--481a4781-a4a9-473f-9ace-fc16e498bc27",
    "attack-pattern--e2088636-a550-43d4-930c-603a33a024e1",
    "attack-pattern--c46f18fa-cad2-4d24-8e9a-aa38bb49ed0d",
    "attack-pattern--6d5ae8f7-0011-49bb-be21-6875ab3e0826",
    "attack-pattern--689f398a-d130-4b8e-b359-5753973d8a8e",
    "attack-pattern--92c88822-861d-4b04-b68d-56cd34ba35e1",
    "attack-pattern--78b26a25-f4e6-48b4-9ce9-aab8bce5db7e",
    "attack-pattern--02cec77c-40b1-418f-920b-a3f3217b669d",
    "attack-pattern--b7dfb4dd-5325-4485-9b63-d57a63bd62d5",
    "attack-pattern--3d802c02-cc2a-4f12-8ec5-2cfd47fa5592",
    "attack-pattern--016a388c-d340-44a8-ab69-26a115e5f9a4",
    "attack-pattern--91582fc7-5103-4d51-b520-1863eac391e1",
    "attack-pattern--7b189c01-b058-4340-8500-e90b3220354a",
    "attack-pattern--3d0968cd-f067-4b8b-a451-df6c521d3dab",
    "attack-pattern--f81f0a08-3399-4006-9175-71a7dca2a917",
    "attack-pattern--ea45586f-e473-43ca-a2b5-44a424eee5e6",
    "attack-pattern--e6771085-d926-419c-831c-64cca5258442",
    "attack-pattern--02efc53c-0695-44f1-8951-b3644058bb60",
    "attack-pattern--3d1fb699-0cd0-420f-82ff-0f6947a10ecc",
    "attack-pattern--9aedbc21-6bff-4826-a8a2-f8a414478fc0",
    "attack-pattern--a562f5cf-fa56-42cb-8bca-3f39c9c8277e"
)

for ($i = 0; $i -lt $newStixIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newStixIds[$i]
}

# --- 2) Copy-edit fixes to technique descriptions (column D) ---

# Row 9: AML.T0043.002 "Black-Box Transfer" description
$ws.Cells.Item(9, 4).Value = "In Black-Box Transfer attacks, the adversary uses one or more proxy models (trained via [Create Proxy ML Model](/techniques/AML.T0005) or [Train Proxy via Replication](/techniques/AML.T0005.001)) they have full access to and are representative of the target model.`nThe adversary uses [White-Box Optimization](/techniques/AML.T0043.000) on the proxy models to generate adversarial examples.`nIf the set of proxy models are close enough to the target model, the adversarial example should generalize from one to another.`nThis means that an attack that works for the proxy models will likely then work for the target model.`nIf the adversary has [ML Model Inference API Access](/techniques/AML.T0040), they may use [Verify Attack](/techniques/AML.T0042) to confirm the attack is working and incorporate that information into their training process.`n"

# Row 83: AML.T0043.000 "White-Box Optimization" description ("manor" -> "manner")
$ws.Cells.Item(83, 4).Value = "In White-Box Optimization, the adversary has full access to the target model and optimizes the adversarial example directly.`nAdversarial examples trained in this manner are most effective against the target model.`n"

# --- 3) Bump created/last modified dates for every technique row (columns F & G) ---
$ws.Range("F2:G83").Value = "12 January 2024"

